$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.885.95"
$ws.Range("E2").Value = "  +5.47%  "
$ws.Range("D3").Value = "2.369.07"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.64"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.38"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "2.368.36"
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.150"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.10"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "2.794.76"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").Value = "60.793.20"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "2.366.26"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.78"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  +9.95%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.06"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("E25").Value = "  +4.41%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.04"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").Value = "  +6.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.76"
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.04"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "0.0₃0737"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").Value = "  +11.58%  "
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.44"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.75"
$ws.Range("E40").Value = "  +11.42%  "
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.27"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.02"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.40"
$ws.Range("E46").Value = "  +7.84%  "
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0216"
$ws.Range("E50").Value = "  +8.11%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("E51").Value = "  +0.96%  "
